# ---------------------------------------------------------------------------
# "Se agregó módulo de reimpresión de acuse individual"
#
# The report sheet gains three new trailing columns (Remesa, Serie Inicial,
# Serie Final, Devuelto push the old W..AC block out to Y..AE) and the
# "Correo" / "ListaParaEnviar" columns are dropped from the header in favour
# of the new ones. Net result: header row grows from A:AC to A:AE.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Widen the merged banner W6:AC6 -> W6:AE6 ---------------------------
$ws.Range("W6:AC6").UnMerge()
$ws.Range("W6:AE6").Merge()

# --- 2. New column widths for W..AE (the tail of the table got narrower,
#        3 fresh columns were appended at 26.33 chars wide) ----------------
$ws.Columns.Item(23).ColumnWidth = 13.498697916666666   # W
$ws.Columns.Item(24).ColumnWidth = 18.498697916666668   # X
$ws.Columns.Item(25).ColumnWidth = 9.998697916666666    # Y
$ws.Columns.Item(26).ColumnWidth = 9.830729166666666    # Z
$ws.Columns.Item(27).ColumnWidth = 8.998697916666666    # AA
$ws.Columns.Item(28).ColumnWidth = 8.998697916666666    # AB
$ws.Columns.Item(29).ColumnWidth = 25.498697916666668   # AC
$ws.Columns.Item(30).ColumnWidth = 25.498697916666668   # AD (new)
$ws.Columns.Item(31).ColumnWidth = 25.498697916666668   # AE (new)

# --- 3. Carry the existing row formatting (rows 1-9) from column AC into
#        the two freshly introduced columns AD and AE ----------------------
for ($r = 1; $r -le 9; $r++) {
  $ws.Cells.Item($r, 29).Copy()
  $ws.Cells.Item($r, 30).PasteSpecial(-4122)
  $ws.Cells.Item($r, 29).Copy()
  $ws.Cells.Item($r, 31).PasteSpecial(-4122)
}

# --- 4. Re-label the header row. Columns B..V keep their meaning; from W
#        onward the layout shifts: Correo is dropped, ListaParaEnviar is
#        replaced, and four report columns are interleaved /appended. -------
$ws.Range("W10").Value = "Estatus"
$ws.Range("Y10").Value = "Remesa"
$ws.Range("Z10").Value = "Serie Inicial"
$ws.Range("AA10").Value = "Serie Final"
$ws.Range("AB10").Value = "Devuelto"
$ws.Range("AC10").Value = "Capturo"
$ws.Range("AD10").Value = "Actualizo"
$ws.Range("AE10").Value = "Enlace"

# Header cells use the bold/shaded style (same one already on B10:AC10) -
# make sure the two new cells AD10/AE10 pick it up too.
$ws.Range("AC10").Copy()
$ws.Range("AD10").PasteSpecial(-4122)
$ws.Range("AC10").Copy()
$ws.Range("AE10").PasteSpecial(-4122)
$ws.Range("AC10").Value = "Capturo"

# --- 5. Printable area now spans through column AE -------------------------
$ws.PageSetup.PrintArea = "`$A`$2:`$AE`$48"

# --- 6. Restore the cursor to the cell it was left on when the workbook was
#        last saved. ---------------------------------------------------------
$ws.Range("A10").Select()
